$d = $word.ActiveDocument
$pStart = $d.Paragraphs.Item(3)
$pEnd = $d.Paragraphs.Item(8)
$target = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t>DML</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> mode</w:t></w:r><w:r><w:t xml:space="preserve"> for INSERT. Finish errors, run</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Then</w:t></w:r><w:r><w:t xml:space="preserve">: give it .net and python dml INSERT for guidelines, convert DML UPDATE and </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>we got data</w:t></w:r><w:r><w:t>!</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>cleanups:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Test the whole ‘data table that was empty</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>’ ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> bulk-insert into it</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>when not in full DML mode, but there should still be some text about data about to change</w:t></w:r><w:r><w:t>1</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xml)
